# Adds the "SwateTemplateMetadata" sheet to the workbook, matching the
# SWATE_templates tooling convention ("Add MetadataSheet to all Templates").
#
# 1) Renames the original data sheet to the template's own name.
# 2) Appends a new "SwateTemplateMetadata" sheet carrying Id/Name/Version/...
#    metadata about the template, styled like the SWATE house style
#    (green/red palette, white bold labels, thick left rule).
# 3) Adds a threaded comment on the Id cell explaining the field and
#    recording the generated id as a reply - mirroring the authoring tool.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename the first (data) worksheet
# ---------------------------------------------------------------------------
$dataSheet = $wb.Worksheets.Item(1)
$dataSheet.Name = "3ASY02_ProteomicsMassSpec"

# ---------------------------------------------------------------------------
# 2) Add the metadata worksheet right after the data sheet
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Add($null, $dataSheet)
$meta.Name = "SwateTemplateMetadata"

# ---------------------------------------------------------------------------
# Colour palette (BGR-packed ints, as consumed by Interior.Color / Font.Color
# / Borders(..).Color -- same encoding RGB() would produce)
# ---------------------------------------------------------------------------
$colWhite      = 16119285   # F5F5F5 - label font colour
$colGreenDark  = 4616993    # 217346 - label fill / section-header border fill
$colGreenLight = 11783583   # 9FCDB3 - value fill
$colGreenDeep  = 3103758    # 0E5C2F - section-header fill
$colRed        = 3809218    # C21F3A - id value fill
$colRule       = 16711422   # FEFEFE - hairline rule colour used everywhere

$xlThin  = 2
$xlThick = 4
$xlEdgeLeft   = 7
$xlEdgeTop    = 8
$xlEdgeBottom = 9
$xlEdgeRight  = 10
$xlTop = -4160

function Set-RuleBorder($rng, $edge, $weight) {
    $b = $rng.Borders.Item($edge)
    $b.Color = $colRule
    $b.Weight = $weight
    $b.LineStyle = 1
}

# Column-A "label" cells: bold off-white text on green, vertical-top aligned,
# thin rule on the left (plus extra top/bottom rule for the first/last row).
function Set-LabelStyle($rng, $fillColor, $topRule, $bottomRule) {
    $rng.Font.Bold = $true
    $rng.Font.Color = $colWhite
    $rng.Interior.Color = $fillColor
    $rng.VerticalAlignment = $xlTop
    Set-RuleBorder $rng $xlEdgeLeft $xlThin
    if ($topRule)    { Set-RuleBorder $rng $xlEdgeTop    $xlThin }
    if ($bottomRule) { Set-RuleBorder $rng $xlEdgeBottom $xlThin }
}

# Column B-F "value" cells: default font, wrap text, vertical-top aligned,
# thick rule on the left + thin rule on the right (plus top/bottom where the
# row is a section boundary).
function Set-ValueStyle($rng, $fillColor, $topRule, $bottomRule) {
    $rng.Interior.Color = $fillColor
    $rng.VerticalAlignment = $xlTop
    $rng.WrapText = $true
    Set-RuleBorder $rng $xlEdgeLeft  $xlThick
    Set-RuleBorder $rng $xlEdgeRight $xlThin
    if ($topRule)    { Set-RuleBorder $rng $xlEdgeTop    $xlThin }
    if ($bottomRule) { Set-RuleBorder $rng $xlEdgeBottom $xlThin }
}

Write-Output "Sheets so far:"
foreach ($s in $wb.Worksheets) { Write-Output $s.Name }

# ---------------------------------------------------------------------------
# 3) Populate the metadata sheet
# ---------------------------------------------------------------------------
# Column A = field label, Column B (and C:F for the multi-value "Tags" row)
# = field value. Rows flagged "section" get the darker header treatment.

$templateId      = "64edd0b7-c2da-4bd0-b6f3-3d150a6151a8"
$templateName    = "Proteomics MassSpec Assay"
$templateVersion = "1.1.3"
$templateDesc    = "This protocol focuses on the measurement of the mass spectrometer, its settings and all other relevant data related to this."
$templateDocs    = "https://github.com/nfdi4plants/SWATE_templates/wiki/3ASY02_ProteomicsMeasurement"
$templateTable   = "annotationTableAverageGoose75"

$meta.Range("A1").Value = "Id"
$meta.Range("B1").Value = $templateId

$meta.Range("A2").Value = "Name"
$meta.Range("B2").Value = $templateName

$meta.Range("A3").Value = "Version"
$meta.Range("B3").Value = "'" + $templateVersion

$meta.Range("A4").Value = "Description"
$meta.Range("B4").Value = $templateDesc

$meta.Range("A5").Value = "Docslink"
$meta.Range("B5").Value = $templateDocs

$meta.Range("A6").Value = "Organisation"

$meta.Range("A7").Value = "Table"
$meta.Range("B7").Value = $templateTable

$meta.Range("A8").Value = "#ER list"

$meta.Range("A9").Value = "ER"
$meta.Range("B9").Value = "PRIDE"

$meta.Range("A10").Value = "ER Term Accession Number"
$meta.Range("A11").Value = "ER Term Source REF"

$meta.Range("A12").Value = "#TAGS list"

$meta.Range("A13").Value = "Tags"
$meta.Range("B13").Value = "Assay"
$meta.Range("C13").Value = "Proteomics"
$meta.Range("D13").Value = "Measurement"
$meta.Range("E13").Value = "Mass spectrometry"
$meta.Range("F13").Value = "MS"

$meta.Range("A14").Value = "Tags Term Accession Number"
$meta.Range("A15").Value = "Tags Term Source REF"

$meta.Range("A16").Value = "#AUTHORS list"

$meta.Range("A17").Value = "Authors Last Name"
$meta.Range("B17").Value = "Maus"

$meta.Range("A18").Value = "Authors First Name"
$meta.Range("B18").Value = "Oliver"

$meta.Range("A19").Value = "Authors Mid Initials"
$meta.Range("A20").Value = "Authors Email"
$meta.Range("A21").Value = "Authors Phone"
$meta.Range("A22").Value = "Authors Fax"
$meta.Range("A23").Value = "Authors Address"
$meta.Range("A24").Value = "Authors Affiliation"

$meta.Range("A25").Value = "#AUTHORS ROLES list"

$meta.Range("A26").Value = "Authors Roles"
$meta.Range("A27").Value = "Authors Roles Term Accession Number"
$meta.Range("A28").Value = "Authors Roles Term Source REF"

Write-Output "Values written"

# ---------------------------------------------------------------------------
# 4) Style the metadata sheet to match the SWATE house style
# ---------------------------------------------------------------------------

# -- Column A labels --------------------------------------------------------
Set-LabelStyle $meta.Range("A1")        $colGreenDark $true  $false
Set-LabelStyle $meta.Range("A2:A7")     $colGreenDark $false $false
Set-LabelStyle $meta.Range("A8")        $colGreenDeep $false $false
Set-LabelStyle $meta.Range("A9:A11")    $colGreenDark $false $false
Set-LabelStyle $meta.Range("A12")       $colGreenDeep $false $false
Set-LabelStyle $meta.Range("A13:A15")   $colGreenDark $false $false
Set-LabelStyle $meta.Range("A16")       $colGreenDeep $false $false
Set-LabelStyle $meta.Range("A17:A24")   $colGreenDark $false $false
Set-LabelStyle $meta.Range("A25")       $colGreenDeep $false $false
Set-LabelStyle $meta.Range("A26:A27")   $colGreenDark $false $false
Set-LabelStyle $meta.Range("A28")       $colGreenDark $false $true

# -- Column B values ----------------------------------------------------
Set-ValueStyle $meta.Range("B1")        $colRed        $true  $false
Set-ValueStyle $meta.Range("B2")        $colGreenLight $false $false
Set-ValueStyle $meta.Range("B3")        $colGreenLight $false $false
Set-ValueStyle $meta.Range("B4")        $colGreenLight $false $false
Set-ValueStyle $meta.Range("B5")        $colGreenLight $false $false
Set-ValueStyle $meta.Range("B6")        $colGreenLight $false $false
Set-ValueStyle $meta.Range("B7")        $colGreenLight $false $false
Set-ValueStyle $meta.Range("B8")        $colGreenDeep  $true  $true
Set-ValueStyle $meta.Range("B9")        $colGreenLight $false $false
Set-ValueStyle $meta.Range("B10")       $colGreenLight $false $false
Set-ValueStyle $meta.Range("B11")       $colGreenLight $false $false
Set-ValueStyle $meta.Range("B12")       $colGreenDeep  $true  $true
Set-ValueStyle $meta.Range("B13:F13")   $colGreenLight $false $false
Set-ValueStyle $meta.Range("B14:F14")   $colGreenLight $false $false
Set-ValueStyle $meta.Range("B15:F15")   $colGreenLight $false $false
Set-ValueStyle $meta.Range("B16")       $colGreenDeep  $true  $true
Set-ValueStyle $meta.Range("B17")       $colGreenLight $false $false
Set-ValueStyle $meta.Range("B18")       $colGreenLight $false $false
Set-ValueStyle $meta.Range("B19")       $colGreenLight $false $false
Set-ValueStyle $meta.Range("B20")       $colGreenLight $false $false
Set-ValueStyle $meta.Range("B21")       $colGreenLight $false $false
Set-ValueStyle $meta.Range("B22")       $colGreenLight $false $false
Set-ValueStyle $meta.Range("B23")       $colGreenLight $false $false
Set-ValueStyle $meta.Range("B24")       $colGreenLight $false $false
Set-ValueStyle $meta.Range("B25")       $colGreenDeep  $true  $true
Set-ValueStyle $meta.Range("B26")       $colGreenLight $false $false
Set-ValueStyle $meta.Range("B27")       $colGreenLight $false $false
Set-ValueStyle $meta.Range("B28")       $colGreenLight $true  $true

Write-Output "Styles applied"
